# Auto-generated edit script: refreshes the Price (D) and Volume(1h) (E)
# columns on the cryptos worksheet, matching the GitHub Actions data-refresh commit.
# D-column price strings are forced to Text format before assignment so that
# values such as '593.67' or '2.931.48' are preserved exactly as text, not coerced
# to numbers (the source data always stores these as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "61.211.84"
$ws.Range("E2").Value = "  +0.45%  "

$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "2.935.09"
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("E4").Value = "  -0.07%  "

$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "593.67"
$ws.Range("E5").Value = "  +0.44%  "

$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "145.24"
$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("E7").Value = "  -0.02%  "

$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.505"
$ws.Range("E8").Value = "  -0.09%  "

$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "6.99"
$ws.Range("E9").Value = "  +4.17%  "

$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("E11").Value = "  -0.47%  "

$ws.Range("E12").Value = "  -0.30%  "

$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "33.83"

$ws.Range("E14").Value = "  +0.13%  "

$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "3.421.09"
$ws.Range("E15").Value = "  +0.79%  "

$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = "61.233.83"
$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("E17").Value = "  +0.03%  "

$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "2.941.00"
$ws.Range("E18").Value = "  +1.20%  "

$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "434.18"
$ws.Range("E19").Value = "  +0.82%  "

$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "13.51"
$ws.Range("E20").Value = "  -0.29%  "

$ws.Range("E21").Value = "  -0.27%  "

$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "7.12"
$ws.Range("E22").Value = "  +0.16%  "

$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "81.88"
$ws.Range("E23").Value = "  +0.65%  "

$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "11.05"
$ws.Range("E24").Value = "  +1.97%  "

$ws.Range("E25").Value = "  -1.23%  "

$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "11.91"
$ws.Range("E26").Value = "  -0.61%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("E28").Value = "  -3.52%  "

$ws.Range("E29").Value = "  -0.40%  "

$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "6.99"
$ws.Range("E30").Value = "  -1.28%  "

$ws.Range("E31").Value = "  +3.35%  "

$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "26.78"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("E33").Value = "  -0.08%  "

$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "0.0₃0874"
$ws.Range("E34").Value = "  +2.09%  "

$ws.Range("E35").Value = "  +0.62%  "

$ws.Range("E36").Value = "  +0.60%  "

$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "2.99"

$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("E39").Value = "  -0.60%  "

$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "8.63"
$ws.Range("E40").Value = "  -0.02%  "

$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "41.67"
$ws.Range("E41").Value = "  +3.30%  "

$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.285"
$ws.Range("E42").Value = "  -2.91%  "

$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "376.74"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("E44").Value = "  -0.46%  "

$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "2.706.95"
$ws.Range("E45").Value = "  +0.08%  "

$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "133.56"
$ws.Range("E46").Value = "  +2.83%  "

$ws.Range("E47").Value = "  -0.06%  "

$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "23.99"
$ws.Range("E48").Value = "  -0.73%  "

$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("E51").Value = "  +0.37%  "

